$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections in column C (English translations) ---
$ws.Range("C5").Value  = "Create"
$ws.Range("C8").Value  = "Created at"
$ws.Range("C14").Value = "Please read and check the service agreement'"
$ws.Range("C15").Value = "Confirm"

# --- New formatting introduced on column D (red font + wrap text) ---
# (touch a single anchor cell so the new font/style is registered and the
#  used range grows to column D, without stamping every row with a style)
$ws.Range("D1").Font.Color = 255
$ws.Range("D1").WrapText = $true

# --- Page setup additions ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection moved (as recorded in the saved view state) ---
$ws.Range("C19").Select()
